$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently has two stacked header rows (row 1 + row 2) above the
# data rows. Row 1 carries the merged-looking fragment headers, row 2 the
# unit headers. Delete row 1 so the unit-header row (and everything below
# it) shifts up by one -- this also trims the trailing blank row 34 for
# free, since the sheet is one row shorter afterwards.
$ws.Rows(1).Delete()

# Row 1 is now the old unit-header row (m3/s / MW / MW / GWh / GWh / GWh).
# Clear it out entirely so we can rebuild the real header row from scratch.
$ws.Range("A1:K1").Clear()

# Rebuild header row 1 with the new column headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 get a distinct style: same font as the data-label cells (font id 2
# in the original styles -- Arial 9) but with no explicit number format
# applied. Achieve this by registering a transient named cell style, using
# it, then deleting the named style again -- the underlying cellXf survives
# (re-parented to xfId 0) while the temporary cellStyle/cellStyleXf entries
# are cleaned back up.
$tempStyle = $wb.Styles.Add("__TempHeaderStyle")
$tempStyle.Font.Name = "Arial"
$tempStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "__TempHeaderStyle"
$wb.Styles.Item("__TempHeaderStyle").Delete() | Out-Null

# Selection moves to A2:K2 in the edited workbook.
$ws.Range("A2:K2").Select() | Out-Null
